$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.198088645935059
$ws.Range("B1").Value = 2.393446922302246
$ws.Range("C1").Value = 4.438894271850586
$ws.Range("D1").Value = 2.715698480606079
$ws.Range("E1").Value = 1.108286619186401
